$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (bold, centered, bordered) from A10 into A11 before writing values
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 694.8099999999999
$ws.Range("C11").Value = 207.72
$ws.Range("D11").Value = 35.93
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 223.53
$ws.Range("G11").Value = 837.84
$ws.Range("H11").Value = 93.78
$ws.Range("I11").Value = 805.86
$ws.Range("J11").Value = 99.37
$ws.Range("K11").Value = 14440.43
$ws.Range("L11").Value = 105.42
$ws.Range("M11").Value = 16.07
$ws.Range("N11").Value = 1.88
$ws.Range("O11").Value = 247.1
$ws.Range("P11").Value = 167.32
$ws.Range("Q11").Value = 13.66
$ws.Range("R11").Value = 21.54
$ws.Range("S11").Value = 479.65
$ws.Range("T11").Value = 44.57
$ws.Range("U11").Value = 2739.75
$ws.Range("V11").Value = ""
$ws.Range("W11").Value = 65.25
$ws.Range("X11").Value = 130.56
$ws.Range("Y11").Value = 105.55
$ws.Range("Z11").Value = 1088.89
$ws.Range("AA11").Value = 198.6
$ws.Range("AB11").Value = 89.89
$ws.Range("AC11").Value = 37.46
$ws.Range("AD11").Value = 231.12
$ws.Range("AE11").Value = 225.65
$ws.Range("AF11").Value = 2581.4
$ws.Range("AG11").Value = 975.29
$ws.Range("AH11").Value = 249.77
$ws.Range("AI11").Value = 142.76
$ws.Range("AJ11").Value = 35.82
$ws.Range("AK11").Value = 397.66
$ws.Range("AL11").Value = 160.63
$ws.Range("AM11").Value = 306.37
$ws.Range("AN11").Value = 7.54
$ws.Range("AO11").Value = 370.87
$ws.Range("AP11").Value = 190.71
$ws.Range("AQ11").Value = 11

Write-Output "row 11 written"